$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new D (Price) and E (Volume(1h)) values.
# Only rows that changed in the source diff are listed; D is optional (some rows
# only had the E/Volume value change).
$updates = @(
    @{ Row = 2;  D = "26.281.76"; E = "  +0.17%  " },
    @{ Row = 3;  D = "1.607.29";  E = "  +0.12%  " },
    @{ Row = 4;  E = "  -0.11%  " },
    @{ Row = 5;  D = "212.95";    E = "  +0.07%  " },
    @{ Row = 6;  E = "  -0.09%  " },
    @{ Row = 7;  D = "0.486";     E = "  +0.35%  " },
    @{ Row = 8;  E = "  +0.03%  " },
    @{ Row = 9;  E = "  -0.50%  " },
    @{ Row = 10; D = "18.33";     E = "  +1.31%  " },
    @{ Row = 11; E = "  -0.36%  " },
    @{ Row = 12; D = "1.832.19";  E = "  +0.29%  " },
    @{ Row = 13; D = "1.606.12";  E = "  +0.11%  " },
    @{ Row = 14; E = "  +0.59%  " },
    @{ Row = 15; E = "  +0.83%  " },
    @{ Row = 16; D = "26.274.01"; E = "  +0.29%  " },
    @{ Row = 17; D = "61.62";     E = "  +1.67%  " },
    @{ Row = 18; E = "  +0.35%  " },
    @{ Row = 19; E = "  -0.21%  " },
    @{ Row = 20; D = "203.97";    E = "  +2.92%  " },
    @{ Row = 21; E = "  +1.20%  " },
    @{ Row = 22; D = "9.31";      E = "  -1.03%  " },
    @{ Row = 23; D = "6.01";      E = "  -0.15%  " },
    @{ Row = 24; D = "1.91";      E = "  +8.95%  " },
    @{ Row = 25; D = "144.23";    E = "  +1.03%  " },
    @{ Row = 26; E = "  +0.01%  " },
    @{ Row = 27; E = "  -6.37%  " },
    @{ Row = 28; D = "15.22";     E = "  +0.29%  " },
    @{ Row = 29; D = "6.56";      E = "  +1.44%  " },
    @{ Row = 30; E = "  +3.66%  " },
    @{ Row = 31; D = "1.18";      E = "  +0.18%  " },
    @{ Row = 32; D = "3.20";      E = "  +2.13%  " },
    @{ Row = 33; D = "2.94";      E = "  -2.15%  " },
    @{ Row = 34; E = "  +3.25%  " },
    @{ Row = 35; E = "  +0.01%  " },
    @{ Row = 36; D = "1.157.16";  E = "  +4.37%  " },
    @{ Row = 37; E = "  +8.88%  " },
    @{ Row = 38; E = "  -0.07%  " },
    @{ Row = 39; E = "  +1.23%  " },
    @{ Row = 40; E = "  +0.31%  " },
    @{ Row = 41; E = "  +0.25%  " },
    @{ Row = 42; D = "0.783";     E = "  +0.58%  " },
    @{ Row = 43; E = "  +2.89%  " },
    @{ Row = 44; D = "1.745.82";  E = "  +0.24%  " },
    @{ Row = 45; D = "91.79";     E = "  -0.78%  " },
    @{ Row = 46; E = "  -1.11%  " },
    @{ Row = 47; D = "54.28";     E = "  +1.36%  " },
    @{ Row = 48; E = "  -0.26%  " },
    @{ Row = 49; E = "  -6.44%  " },
    @{ Row = 50; E = "  -0.63%  " },
    @{ Row = 51; E = "  -0.36%  " }
)

# The "Price" column (D) holds values such as "26.281.76" or "212.95" that are
# formatted as plain text in the source file (t="inlineStr"), not numbers -
# some of them even use multiple "." characters as thousands separators, so
# they can never be read back as a single numeric value. Assigning a
# single-dot numeric-looking string straight to Range.Value (e.g. "212.95")
# would make Excel auto-convert it to a real floating point number, which
# does not match the original text-based layout.
#
# To force these to stay plain text *without* changing the cell's number
# format/style (the original cells use the default style, no explicit
# NumberFormat), we stage the text in a helper cell that has an explicit
# Text ("@") number format, copy it, and paste-special *values only* into
# the destination. PasteSpecial(values) carries over the text-ness of the
# source without carrying over its number format/style, so the destination
# cell keeps its original (default) style while its content becomes text.
$helper = $ws.Cells.Item(200, 26)   # far-away scratch cell (Z200)
$helper.NumberFormat = "@"

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $helper.Value = $u.D
        $helper.Copy()
        $ws.Cells.Item($r, 4).PasteSpecial(-4163)  # xlPasteValues
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}

$helper.Clear()
$excel.CutCopyMode = $false
